$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp note (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Mayo de 2020 a las 14:05"

# --- Reinsert "Azerbaiyan" earlier in the ranking (row 69), pushing
#     Luxemburgo and Irak down one row each with their previous data,
#     and giving Azerbaiyan its updated figures. ---

# Row 69: now Azerbaiyan, with new totals
$ws.Range("A69").Value = "Azerbaiyan"
$ws.Range("B69").Value = 3982
$ws.Range("C69").Value = 127
$ws.Range("D69").Value = 2506
$ws.Range("E69").Value = 1427
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 3
$ws.Range("H69").Value = 49

# Row 70: now Luxemburgo (keeps its previous totals)
$ws.Range("A70").Value = "Luxemburgo"
$ws.Range("B70").Value = 3981
$ws.Range("C70").Value = 0
$ws.Range("D70").Value = 3748
$ws.Range("E70").Value = 124
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 109

# Row 71: now Irak (keeps its previous totals)
$ws.Range("A71").Value = "Irak"
$ws.Range("B71").Value = 3964
$ws.Range("C71").Value = 0
$ws.Range("D71").Value = 2532
$ws.Range("E71").Value = 1285
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 147

# --- Update covid figures for several other countries (rows unaffected
#     by reordering) ---

# Row 27: Suecia
$ws.Range("B27").Value = 33188
$ws.Range("C27").Value = 379
$ws.Range("E27").Value = 24225
$ws.Range("G27").Value = 67
$ws.Range("H27").Value = 3992

# Row 47: Dinamarca
$ws.Range("B47").Value = 11289
$ws.Range("C47").Value = 59
$ws.Range("D47").Value = 9836
$ws.Range("E47").Value = 892

# Row 88: Republica de Macedonia
$ws.Range("B88").Value = 1941
$ws.Range("C88").Value = 20
$ws.Range("D88").Value = 1411
$ws.Range("E88").Value = 417
$ws.Range("G88").Value = 1
$ws.Range("H88").Value = 113

# Row 104: Libano
$ws.Range("B104").Value = 1097
$ws.Range("C104").Value = 11
$ws.Range("D104").Value = 667
$ws.Range("E104").Value = 404

# Row 105: Sri Lanka
$ws.Range("B105").Value = 1078
$ws.Range("C105").Value = 10
$ws.Range("E105").Value = 409
